# Generate Report for Archive
#
# Semantic change (per the target diff):
#  - Status "Ready for handoff" becomes "In Translation" for the two files
#    that are still mid-translation (27d94dc7..., a377dd55...).
#  - The row order of a377dd55... and 52870d72... swaps (a377dd55 now sorts
#    before 52870d72) on every sheet (Overview, zh-cn, de-de), carrying each
#    file's own data + hyperlink along with it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 3 (27d94dc7...): status columns E/F flip to "In Translation"
$ov.Range("E3").Value = "In Translation"
$ov.Range("F3").Value = "In Translation"

# Row 4 now holds a377dd55... (was 52870d72...)
$ov.Range("A4").Value = "a377dd55-3b7f-40ad-9d1d-faf597ac9235.md"
$ov.Range("B4").Value = "e2e\a377dd55-3b7f-40ad-9d1d-faf597ac9235.md"
$ov.Range("E4").Value = "In Translation"
$ov.Range("F4").Value = "In Translation"
$ov.Range("G4").Value = "2016-08-18 20:41:47"

# Row 5 now holds 52870d72... (was a377dd55...)
$ov.Range("A5").Value = "52870d72-cb67-4dde-b5bb-6a6e6d11f791.md"
$ov.Range("B5").Value = "e2e\52870d72-cb67-4dde-b5bb-6a6e6d11f791.md"
$ov.Range("E5").Value = "Ready for handoff"
$ov.Range("F5").Value = "Ready for handoff"
$ov.Range("G5").Value = "2016-08-18 20:40:34"

# Rebuild hyperlinks in B2:B5 so column B keeps matching its row's file.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb0c854b7ad8c3b7cec9bed44ec340174e2e39fb/e2e/8603798e-08d2-4224-b200-d3d69149f632.md", "", "", "e2e\8603798e-08d2-4224-b200-d3d69149f632.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77ea68f9a4a36880f655d0968aa545fb2651e3e/e2e/27d94dc7-3b87-455c-9825-37e01d184ffb.md", "", "", "e2e\27d94dc7-3b87-455c-9825-37e01d184ffb.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77ea68f9a4a36880f655d0968aa545fb2651e3e/e2e/a377dd55-3b7f-40ad-9d1d-faf597ac9235.md", "", "", "e2e\a377dd55-3b7f-40ad-9d1d-faf597ac9235.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d5449d797d7662fe7f86abce130cc70f3a73b8/e2e/52870d72-cb67-4dde-b5bb-6a6e6d11f791.md", "", "", "e2e\52870d72-cb67-4dde-b5bb-6a6e6d11f791.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 3 (27d94dc7...): Status -> "In Translation"
$zh.Range("C3").Value = "In Translation"

# Row 4 now holds a377dd55... (was 52870d72...)
$zh.Range("A4").Value = "a377dd55-3b7f-40ad-9d1d-faf597ac9235.md"
$zh.Range("C4").Value = "In Translation"
$zh.Range("G4").Value = "a377dd55-3b7f-40ad-9d1d-faf597ac9235.9de346615b3b055c40d3c409e8c8dcc32e072631.zh-cn.xlf"
$zh.Range("H4").Value = "2016-08-18 20:41:42"

# Row 5 now holds 52870d72... (was a377dd55...)
$zh.Range("A5").Value = "52870d72-cb67-4dde-b5bb-6a6e6d11f791.md"
$zh.Range("C5").Value = "Ready for handoff"
$zh.Range("G5").Value = "52870d72-cb67-4dde-b5bb-6a6e6d11f791.2e149822f8326a8b35d7d4bf6c1d37e51e259880.zh-cn.xlf"
$zh.Range("H5").Value = "2016-08-18 20:40:29"

# Rebuild hyperlinks: A2:A5 per-row file link, I2 keeps its own language-repo link.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb0c854b7ad8c3b7cec9bed44ec340174e2e39fb/e2e/8603798e-08d2-4224-b200-d3d69149f632.md", "", "", "8603798e-08d2-4224-b200-d3d69149f632.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a1390d17640f0bc2c05e4c7bf35216914fa5881e/e2e/8603798e-08d2-4224-b200-d3d69149f632.md", "", "", "8603798e-08d2-4224-b200-d3d69149f632.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77ea68f9a4a36880f655d0968aa545fb2651e3e/e2e/27d94dc7-3b87-455c-9825-37e01d184ffb.md", "", "", "27d94dc7-3b87-455c-9825-37e01d184ffb.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77ea68f9a4a36880f655d0968aa545fb2651e3e/e2e/a377dd55-3b7f-40ad-9d1d-faf597ac9235.md", "", "", "a377dd55-3b7f-40ad-9d1d-faf597ac9235.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d5449d797d7662fe7f86abce130cc70f3a73b8/e2e/52870d72-cb67-4dde-b5bb-6a6e6d11f791.md", "", "", "52870d72-cb67-4dde-b5bb-6a6e6d11f791.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 3 (27d94dc7...): Status -> "In Translation"
$de.Range("C3").Value = "In Translation"

# Row 4 now holds a377dd55... (was 52870d72...)
$de.Range("A4").Value = "a377dd55-3b7f-40ad-9d1d-faf597ac9235.md"
$de.Range("C4").Value = "In Translation"
$de.Range("G4").Value = "a377dd55-3b7f-40ad-9d1d-faf597ac9235.9de346615b3b055c40d3c409e8c8dcc32e072631.de-de.xlf"
$de.Range("H4").Value = "2016-08-18 20:41:47"

# Row 5 now holds 52870d72... (was a377dd55...)
$de.Range("A5").Value = "52870d72-cb67-4dde-b5bb-6a6e6d11f791.md"
$de.Range("C5").Value = "Ready for handoff"
$de.Range("G5").Value = "52870d72-cb67-4dde-b5bb-6a6e6d11f791.2e149822f8326a8b35d7d4bf6c1d37e51e259880.de-de.xlf"
$de.Range("H5").Value = "2016-08-18 20:40:34"

# Rebuild hyperlinks: A2:A5 per-row file link, I2 keeps its own language-repo link.
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb0c854b7ad8c3b7cec9bed44ec340174e2e39fb/e2e/8603798e-08d2-4224-b200-d3d69149f632.md", "", "", "8603798e-08d2-4224-b200-d3d69149f632.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0791e7d5825725caf395f69061604cd958f314e7/e2e/8603798e-08d2-4224-b200-d3d69149f632.md", "", "", "8603798e-08d2-4224-b200-d3d69149f632.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77ea68f9a4a36880f655d0968aa545fb2651e3e/e2e/27d94dc7-3b87-455c-9825-37e01d184ffb.md", "", "", "27d94dc7-3b87-455c-9825-37e01d184ffb.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77ea68f9a4a36880f655d0968aa545fb2651e3e/e2e/a377dd55-3b7f-40ad-9d1d-faf597ac9235.md", "", "", "a377dd55-3b7f-40ad-9d1d-faf597ac9235.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d5449d797d7662fe7f86abce130cc70f3a73b8/e2e/52870d72-cb67-4dde-b5bb-6a6e6d11f791.md", "", "", "52870d72-cb67-4dde-b5bb-6a6e6d11f791.md") | Out-Null

Write-Output "done"
